# Scheduled-runner update: refresh Maduin Profits crafting-leve profit/loss figures
# across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets (source data refresh; CUL unaffected).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 221
$ws.Range("I12").Value = 204.5
$ws.Range("K12").Value = 204.5
$ws.Range("M12").Value = -34.5
# Row 64
$ws.Range("I64").Value = 20000
$ws.Range("K64").Value = 20000
$ws.Range("M64").Value = -19752
# Row 67
$ws.Range("I67").Value = 20000
$ws.Range("K67").Value = 20000
$ws.Range("M67").Value = -19142
# Row 70
$ws.Range("H70").Value = 2712.125
$ws.Range("I70").Value = 2712.125
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 8136.375
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -7866.375
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 2712.125
$ws.Range("I73").Value = 2712.125
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 8136.375
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -7200.375
$ws.Range("N73").ClearContents()
# Row 74
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
# Row 76
$ws.Range("H76").Value = 4099.7144
$ws.Range("I76").Value = 4259.8
$ws.Range("K76").Value = 4259.8
$ws.Range("M76").Value = -3944.8
# Row 77
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
# Row 79
$ws.Range("H79").Value = 4099.7144
$ws.Range("I79").Value = 4259.8
$ws.Range("K79").Value = 4259.8
$ws.Range("M79").Value = -3167.8
# Row 98
$ws.Range("H98").Value = 1994.5
$ws.Range("I98").Value = 2329.3333
$ws.Range("J98").Value = 990
$ws.Range("K98").Value = 2329.3333
$ws.Range("L98").Value = 990
$ws.Range("M98").Value = -831.3332999999998
$ws.Range("N98").Value = -3986
# Row 113
$ws.Range("H113").Value = 14863.444
$ws.Range("I113").Value = 16252.286
$ws.Range("K113").Value = 16252.286
$ws.Range("M113").Value = -12998.286
# Row 116
$ws.Range("H116").Value = 1281.6666
$ws.Range("J116").Value = 1281.6666
$ws.Range("L116").Value = 1281.6666
$ws.Range("N116").Value = -8165.6666
# Row 122
$ws.Range("H122").Value = 1994.5
$ws.Range("I122").Value = 2329.3333
$ws.Range("J122").Value = 990
$ws.Range("K122").Value = 6987.999899999999
$ws.Range("L122").Value = 2970
$ws.Range("M122").Value = -4537.999899999999
$ws.Range("N122").Value = -7870
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 810.4
$ws.Range("I74").Value = 810.4
$ws.Range("K74").Value = 810.4
$ws.Range("M74").Value = 63.60000000000002
# Row 77
$ws.Range("H77").Value = 810.4
$ws.Range("I77").Value = 810.4
$ws.Range("K77").Value = 4052
$ws.Range("M77").Value = 316
# Row 110
$ws.Range("H110").Value = 1039.6666
$ws.Range("I110").Value = 909.75
$ws.Range("J110").Value = 1299.5
$ws.Range("K110").Value = 909.75
$ws.Range("L110").Value = 1299.5
$ws.Range("M110").Value = 1135.25
$ws.Range("N110").Value = -5389.5
# Row 132
$ws.Range("H132").Value = 1592.75
$ws.Range("J132").Value = 1331
$ws.Range("L132").Value = 3993
$ws.Range("N132").Value = -9053
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3256.2144
$ws.Range("I86").Value = 3012.7144
$ws.Range("J86").Value = 3499.7144
$ws.Range("K86").Value = 3012.7144
$ws.Range("L86").Value = 3499.7144
$ws.Range("M86").Value = -1889.7144
$ws.Range("N86").Value = -5745.7144
# Row 89
$ws.Range("H89").Value = 3256.2144
$ws.Range("I89").Value = 3012.7144
$ws.Range("J89").Value = 3499.7144
$ws.Range("K89").Value = 15063.572
$ws.Range("L89").Value = 17498.572
$ws.Range("M89").Value = -9447.572
$ws.Range("N89").Value = -28730.572
# Row 94
$ws.Range("H94").Value = 2907.125
$ws.Range("I94").Value = 814.5
$ws.Range("K94").Value = 814.5
$ws.Range("M94").Value = -363.5
# Row 105
$ws.Range("H105").Value = 3402.2856
$ws.Range("I105").Value = 3402.2856
$ws.Range("K105").Value = 3402.2856
$ws.Range("M105").Value = -1655.2856
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1557.9333
$ws.Range("I31").Value = 1146.7778
$ws.Range("K31").Value = 1146.7778
$ws.Range("M31").Value = -851.7778000000001
# Row 34
$ws.Range("H34").Value = 1557.9333
$ws.Range("I34").Value = 1146.7778
$ws.Range("K34").Value = 1146.7778
$ws.Range("M34").Value = -944.7778000000001
# Row 122
$ws.Range("H122").Value = 2894.4285
$ws.Range("I122").Value = 2941.8333
$ws.Range("K122").Value = 8825.499899999999
$ws.Range("M122").Value = -6375.499899999999
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 7506
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 7506
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 7506
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -9502
# Row 83
$ws.Range("H83").Value = 7506
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 7506
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 37530
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -47514
# Row 102
$ws.Range("H102").Value = 4343.143
$ws.Range("I102").Value = 3567
$ws.Range("K102").Value = 3567
$ws.Range("M102").Value = -1945
# Row 122
$ws.Range("H122").Value = 2301.2778
$ws.Range("I122").Value = 2150.9375
$ws.Range("K122").Value = 6452.8125
$ws.Range("M122").Value = -4002.8125
$ws = $wb.Worksheets.Item("LTW")
# Row 38
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
# Row 55
$ws.Range("H55").Value = 664.0714
$ws.Range("J55").Value = 799.1818
$ws.Range("L55").Value = 799.1818
$ws.Range("N55").Value = -1145.1818
# Row 93
$ws.Range("H93").Value = 892.25
$ws.Range("I93").Value = 917.8
$ws.Range("K93").Value = 917.8
$ws.Range("M93").Value = 330.2
# Row 136
$ws.Range("H136").Value = 2567.6
$ws.Range("I136").Value = 2491.2942
$ws.Range("K136").Value = 7473.882599999999
$ws.Range("M136").Value = -4923.882599999999
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 5492.5557
$ws.Range("I81").Value = 5738.8335
$ws.Range("K81").Value = 11477.667
$ws.Range("M81").Value = -10416.667
# Row 84
$ws.Range("H84").Value = 5492.5557
$ws.Range("I84").Value = 5738.8335
$ws.Range("K84").Value = 57388.335
$ws.Range("M84").Value = -52084.335
# Row 104
$ws.Range("H104").Value = 23067
$ws.Range("J104").Value = 23067
$ws.Range("L104").Value = 23067
$ws.Range("N104").Value = -30055
# Row 107
$ws.Range("H107").Value = 478.58334
$ws.Range("I107").Value = 482.77777
$ws.Range("K107").Value = 1448.33331
$ws.Range("M107").Value = 471.66669
# Row 122
$ws.Range("H122").Value = 2423.125
$ws.Range("I122").Value = 2191
$ws.Range("K122").Value = 6573
$ws.Range("M122").Value = -4123
# Row 126
$ws.Range("H126").Value = 2021.2
$ws.Range("I126").Value = 2021.2
$ws.Range("K126").Value = 6063.6
$ws.Range("M126").Value = -3593.6
